$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 6250.3335
$ws.Range("I11").Value = 6250.3335
$ws.Range("K11").Value = 6250.3335
$ws.Range("M11").Value = -6110.3335

$ws.Range("H15").Value = 2433.4722
$ws.Range("I15").Value = 2433.4722
$ws.Range("K15").Value = 7300.4166
$ws.Range("M15").Value = -7131.4166

$ws.Range("H19").Value = 970.6667
$ws.Range("J19").Value = 933.5625
$ws.Range("L19").Value = 933.5625
$ws.Range("N19").Value = -1283.5625

$ws.Range("H51").Value = 15714.833
$ws.Range("J51").Value = 6963.3335
$ws.Range("L51").Value = 6963.3335
$ws.Range("N51").Value = -7931.3335

$ws.Range("H87").Value = 53332.668
$ws.Range("J87").Value = 53332.668
$ws.Range("L87").Value = 53332.668
$ws.Range("N87").Value = -55828.668

$ws.Range("H90").Value = 53332.668
$ws.Range("J90").Value = 53332.668
$ws.Range("L90").Value = 159998.004
$ws.Range("N90").Value = -172478.004

$ws.Range("H98").Value = 34486108
$ws.Range("I98").Value = 40003332
$ws.Range("K98").Value = 40003332
$ws.Range("M98").Value = -40001834

$ws.Range("H122").Value = 34486108
$ws.Range("I122").Value = 40003332
$ws.Range("K122").Value = 120009996
$ws.Range("M122").Value = -120007546

$ws.Range("H133").Value = 115385
$ws.Range("J133").Value = 115385
$ws.Range("L133").Value = 115385
$ws.Range("N133").Value = -125505

$ws.Range("H135").Value = 313092.4
$ws.Range("J135").Value = 2555.5
$ws.Range("L135").Value = 22999.5
$ws.Range("N135").Value = -28069.5

$ws.Range("H137").Value = 2354.3684
$ws.Range("I137").Value = 2535.5833
$ws.Range("K137").Value = 7606.749899999999
$ws.Range("M137").Value = -5056.749899999999

$ws.Range("H138").Value = 3855.3408
$ws.Range("I138").Value = 916.05
$ws.Range("J138").Value = 6304.75
$ws.Range("K138").Value = 2748.15
$ws.Range("L138").Value = 18914.25
$ws.Range("M138").Value = 2391.85
$ws.Range("N138").Value = -29194.25

$ws.Range("H141").Value = 1518.6666
$ws.Range("I141").Value = 1518.6666
$ws.Range("K141").Value = 4555.9998
$ws.Range("M141").Value = 624.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2949.3076
$ws.Range("I2").Value = 1578.2632
$ws.Range("J2").Value = 6670.7144
$ws.Range("K2").Value = 1578.2632
$ws.Range("L2").Value = 6670.7144
$ws.Range("M2").Value = -1465.2632
$ws.Range("N2").Value = -6896.7144

$ws.Range("H32").Value = 2233887.5
$ws.Range("I32").Value = 2316493.5
$ws.Range("J32").Value = 3529.5
$ws.Range("K32").Value = 2316493.5
$ws.Range("L32").Value = 3529.5
$ws.Range("M32").Value = -2316206.5
$ws.Range("N32").Value = -4103.5

$ws.Range("H61").Value = 3759.4736
$ws.Range("I61").Value = 2141.2415
$ws.Range("J61").Value = 8973.777
$ws.Range("K61").Value = 2141.2415
$ws.Range("L61").Value = 8973.777
$ws.Range("M61").Value = -1929.2415
$ws.Range("N61").Value = -9397.777

$ws.Range("H88").Value = 1698.0741
$ws.Range("I88").Value = 1370.7
$ws.Range("J88").Value = 1890.6471
$ws.Range("K88").Value = 1370.7
$ws.Range("L88").Value = 1890.6471
$ws.Range("M88").Value = -964.7
$ws.Range("N88").Value = -2702.6471

$ws.Range("H91").Value = 1698.0741
$ws.Range("I91").Value = 1370.7
$ws.Range("J91").Value = 1890.6471
$ws.Range("K91").Value = 1370.7
$ws.Range("L91").Value = 1890.6471
$ws.Range("M91").Value = 33.29999999999995
$ws.Range("N91").Value = -4698.6471

$ws.Range("H97").Value = 8075.4375
$ws.Range("I97").Value = 590.5454999999999
$ws.Range("K97").Value = 590.5454999999999
$ws.Range("M97").Value = -94.54549999999995

$ws.Range("H116").Value = 2949.3076
$ws.Range("I116").Value = 1578.2632
$ws.Range("J116").Value = 6670.7144
$ws.Range("K116").Value = 1578.2632
$ws.Range("L116").Value = 6670.7144
$ws.Range("M116").Value = 715.7367999999999
$ws.Range("N116").Value = -11258.7144

$ws.Range("H122").Value = 27339.3
$ws.Range("I122").Value = 35199
$ws.Range("K122").Value = 105597
$ws.Range("M122").Value = -103147

$ws.Range("H132").Value = 2888.868
$ws.Range("I132").Value = 1673.8857
$ws.Range("K132").Value = 5021.6571
$ws.Range("M132").Value = -2491.6571

$ws.Range("H136").Value = 3759.4736
$ws.Range("I136").Value = 2141.2415
$ws.Range("J136").Value = 8973.777
$ws.Range("K136").Value = 6423.7245
$ws.Range("L136").Value = 26921.331
$ws.Range("M136").Value = -3873.7245
$ws.Range("N136").Value = -32021.331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2949.3076
$ws.Range("I3").Value = 1578.2632
$ws.Range("J3").Value = 6670.7144
$ws.Range("K3").Value = 1578.2632
$ws.Range("L3").Value = 6670.7144
$ws.Range("M3").Value = -1464.2632
$ws.Range("N3").Value = -6898.7144

$ws.Range("H96").Value = 11321.5
$ws.Range("I96").Value = 11321.5
$ws.Range("K96").Value = 11321.5
$ws.Range("M96").Value = -8575.5

$ws.Range("H99").Value = 3457.9412
$ws.Range("I99").Value = 1367.4
$ws.Range("K99").Value = 1367.4
$ws.Range("M99").Value = 130.5999999999999

$ws.Range("H134").Value = 3925.8
$ws.Range("I134").Value = 3036.9272
$ws.Range("K134").Value = 9110.7816
$ws.Range("M134").Value = -6575.7816

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5151.533
$ws.Range("I31").Value = 2931.4443
$ws.Range("J31").Value = 6967.9697
$ws.Range("K31").Value = 2931.4443
$ws.Range("L31").Value = 6967.9697
$ws.Range("M31").Value = -2636.4443
$ws.Range("N31").Value = -7557.9697

$ws.Range("H34").Value = 5151.533
$ws.Range("I34").Value = 2931.4443
$ws.Range("J34").Value = 6967.9697
$ws.Range("K34").Value = 2931.4443
$ws.Range("L34").Value = 6967.9697
$ws.Range("M34").Value = -2729.4443
$ws.Range("N34").Value = -7371.9697

$ws.Range("I58").Value = 15626922
$ws.Range("J58").Value = 6094.577
$ws.Range("K58").Value = 15626922
$ws.Range("L58").Value = 6094.577
$ws.Range("M58").Value = -15626719
$ws.Range("N58").Value = -6500.577

$ws.Range("H132").Value = 3518.4614
$ws.Range("I132").Value = 1593.3077
$ws.Range("K132").Value = 4779.9231
$ws.Range("M132").Value = -2249.9231

$ws.Range("I136").Value = 15626922
$ws.Range("J136").Value = 6094.577
$ws.Range("K136").Value = 46880766
$ws.Range("L136").Value = 18283.731
$ws.Range("M136").Value = -46878216
$ws.Range("N136").Value = -23383.731

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 27777992
$ws.Range("I14").Value = 27777992
$ws.Range("K14").Value = 83333976
$ws.Range("M14").Value = -83333803

$ws.Range("H99").Value = 11333
$ws.Range("I99").Value = 6999.5
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 20998.5
$ws.Range("L99").Value = 60000
$ws.Range("M99").Value = -18752.5
$ws.Range("N99").Value = -64492

$ws.Range("H107").Value = 28571966
$ws.Range("I107").Value = 299
$ws.Range("K107").Value = 897
$ws.Range("M107").Value = 1023

$ws.Range("H122").Value = 2831601.2
$ws.Range("I122").Value = 5658163.5
$ws.Range("J122").Value = 5038.8
$ws.Range("K122").Value = 50923471.5
$ws.Range("L122").Value = 45349.2
$ws.Range("M122").Value = -50921021.5
$ws.Range("N122").Value = -50249.2

$ws.Range("H132").Value = 9176.076999999999
$ws.Range("I132").Value = 7400
$ws.Range("J132").Value = 10286.125
$ws.Range("K132").Value = 66600
$ws.Range("L132").Value = 92575.125
$ws.Range("M132").Value = -64070
$ws.Range("N132").Value = -97635.125

$ws.Range("H141").Value = 8398.556
$ws.Range("I141").Value = 2597.8333
$ws.Range("J141").Value = 20000
$ws.Range("K141").Value = 7793.499899999999
$ws.Range("L141").Value = 60000
$ws.Range("M141").Value = -2613.499899999999
$ws.Range("N141").Value = -70360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5878.423
$ws.Range("I113").Value = 2871.2856
$ws.Range("J113").Value = 6986.316
$ws.Range("K113").Value = 2871.2856
$ws.Range("L113").Value = 6986.316
$ws.Range("M113").Value = -701.2856000000002
$ws.Range("N113").Value = -11326.316

$ws.Range("H122").Value = 43892.617
$ws.Range("I122").Value = 66926.5
$ws.Range("J122").Value = 7038.4
$ws.Range("K122").Value = 200779.5
$ws.Range("L122").Value = 21115.2
$ws.Range("M122").Value = -198329.5
$ws.Range("N122").Value = -26015.2

$ws.Range("H126").Value = 2869
$ws.Range("I126").Value = 2433
$ws.Range("K126").Value = 7299
$ws.Range("M126").Value = -4829

$ws.Range("H132").Value = 2940.8445
$ws.Range("I132").Value = 2836.2163
$ws.Range("J132").Value = 3424.75
$ws.Range("K132").Value = 8508.6489
$ws.Range("L132").Value = 10274.25
$ws.Range("M132").Value = -5978.6489
$ws.Range("N132").Value = -15334.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5684.64
$ws.Range("I7").Value = 3323.0833
$ws.Range("K7").Value = 3323.0833
$ws.Range("M7").Value = -3211.0833

$ws.Range("H22").Value = 3385.9
$ws.Range("J22").Value = 4991.1665
$ws.Range("L22").Value = 4991.1665
$ws.Range("N22").Value = -5581.1665

$ws.Range("H27").Value = 3385.9
$ws.Range("J27").Value = 4991.1665
$ws.Range("L27").Value = 4991.1665
$ws.Range("N27").Value = -5205.1665

$ws.Range("H93").Value = 5247.2573
$ws.Range("I93").Value = 2927.5
$ws.Range("J93").Value = 9173
$ws.Range("K93").Value = 2927.5
$ws.Range("L93").Value = 9173
$ws.Range("M93").Value = -1679.5
$ws.Range("N93").Value = -11669

$ws.Range("H126").Value = 5684.64
$ws.Range("I126").Value = 3323.0833
$ws.Range("K126").Value = 9969.249899999999
$ws.Range("M126").Value = -7499.249899999999

$ws.Range("H136").Value = 8152.5884
$ws.Range("I136").Value = 6931.3335
$ws.Range("J136").Value = 8414.286
$ws.Range("K136").Value = 20794.0005
$ws.Range("L136").Value = 25242.858
$ws.Range("M136").Value = -18244.0005
$ws.Range("N136").Value = -30342.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 125129870
$ws.Range("I136").Value = 333335330
$ws.Range("J136").Value = 206600
$ws.Range("K136").Value = 1000005990
$ws.Range("L136").Value = 619800
$ws.Range("M136").Value = -1000003440
$ws.Range("N136").Value = -624900
